# Adds a third slide ("Gene A" / "Gene B" mini diagram) reproducing the
# layout used on slides 1-2 of the deck (small 1200pt bold maroon labels
# on a white textbox, joined by maroon arrow connectors).

function EMU($v) {
    # PowerPoint COM geometry is expressed in points; the OOXML is EMU
    # (1 pt = 12700 EMU). A tiny epsilon nudges values that would
    # otherwise truncate one EMU short after the float32 round-trip the
    # host applies internally.
    return ([double]$v / 12700.0) + 0.00002
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Duplicate slide 2: this is the easiest way to get a new slide that
#    already carries the same closing <p:extLst>/<p:clrMapOvr> shape as
#    every other slide in the deck, plus a ready-made "TextBox 4"
#    (white fill, no line, centered 1200/bold/maroon run) we can reuse
#    verbatim for the "Gene A" label.
# ---------------------------------------------------------------------
$dup = $p.Slides.Item(2).Duplicate()
$s = $p.Slides.Item(3)

# Keep only the "TextBox 4" shape (Id 5); drop the other 17 clones.
for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Id -ne 5) {
        $sh.Delete()
    }
}

# ---------------------------------------------------------------------
# 2. Turn the kept textbox into the "Gene A" label.
# ---------------------------------------------------------------------
$geneA = $s.Shapes.Item(1)
$geneA.Name = "TextBox 4"
$geneA.Left = (EMU 908392)
$geneA.Top = (EMU 452083)
$geneA.Width = (EMU 1151287)
$geneA.Height = (EMU 276999)
$geneA.TextFrame.TextRange.Text = "Gene A"

# ---------------------------------------------------------------------
# 3. The slide's shape-id counter must land on 20/21/23 for the three
#    remaining shapes (matching the authored file). Burn through the
#    intermediate ids with disposable textboxes.
# ---------------------------------------------------------------------
function BurnIds($n) {
    for ($k = 0; $k -lt $n; $k++) {
        $junk = $s.Shapes.AddTextbox(1, 0, 0, 1, 1)
        $junk.Delete()
    }
}

BurnIds 17
$geneB = $s.Shapes.AddTextbox(1, (EMU 908392), (EMU 1714395), (EMU 1151287), (EMU 276999))
$geneB.Name = "TextBox 19"
$geneB.Fill.ForeColor.RGB = 16777215
$geneB.Line.Visible = $false
$geneB.TextFrame.WordWrap = $true
$geneB.TextFrame.AutoSize = 1
$geneB.Left = (EMU 908392)
$geneB.Top = (EMU 1714395)
$geneB.Width = (EMU 1151287)
$geneB.Height = (EMU 276999)
$tr = $geneB.TextFrame.TextRange
$tr.Text = "Gene B"
$tr.ParagraphFormat.Alignment = 2
$tr.Font.Size = 12
$tr.Font.Bold = $true
$tr.Font.Color.RGB = 128

# ---------------------------------------------------------------------
# 4. Connectors: copy one of slide 1's straight arrow connectors so the
#    shape-style quick-style block (<p:style>) comes along for free,
#    then reposition/rename/retext it.
# ---------------------------------------------------------------------
$cxnSrc = $p.Slides.Item(1).Shapes.Item(4)

$cxnSrc.Copy()
$cxn1 = $s.Shapes.Paste().Item(1)
$cxn1.Name = "Straight Arrow Connector 20"
$cxn1.Left = (EMU 1160342)
$cxn1.Top = (EMU 452083)
$cxn1.Width = (EMU 0)
$cxn1.Height = (EMU 392448)

BurnIds 1
$cxnSrc.Copy()
$cxn2 = $s.Shapes.Paste().Item(1)
$cxn2.Name = "Straight Arrow Connector 22"
$cxn2.Left = (EMU 1160342)
$cxn2.Top = (EMU 1690740)
$cxn2.Width = (EMU 0)
$cxn2.Height = (EMU 392448)

Write-Host "Slide 3 built with" $s.Shapes.Count "shapes"
